# Updates cryptos list: Price (column D) and Volume(1h) (column E) values
# for the rows whose figures moved, per the upstream GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") cells are free-form text (e.g. "1.00", "68.341.27")
# rather than real numbers, so force Text format before writing the new
# string -- otherwise Excel would parse them as numeric and mangle values
# like "1.00" -> 1 or "601.76" -> 601.75999999999999.
$priceRows = @(2,3,4,5,6,7,12,14,15,16,17,18,19,22,23,25,26,27,28,31,32,34,35,36,37,38,39,42,43,46,47,49,50)
foreach ($r in $priceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Cells.Item(2, 4).Value = "68.367.27"
$ws.Cells.Item(2, 5).Value = "  +0.77%  "
$ws.Cells.Item(3, 4).Value = "3.853.26"
$ws.Cells.Item(3, 5).Value = "  +0.10%  "
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$ws.Cells.Item(5, 4).Value = "601.76"
$ws.Cells.Item(5, 5).Value = "  +0.59%  "
$ws.Cells.Item(6, 4).Value = "172.41"
$ws.Cells.Item(6, 5).Value = "  +3.82%  "
$ws.Cells.Item(7, 4).Value = "3.850.72"
$ws.Cells.Item(7, 5).Value = "  +0.05%  "
$ws.Cells.Item(8, 5).Value = "  +0.02%  "
$ws.Cells.Item(9, 5).Value = "  +0.77%  "
$ws.Cells.Item(10, 5).Value = "  +2.01%  "
$ws.Cells.Item(11, 5).Value = "  +3.48%  "
$ws.Cells.Item(12, 4).Value = "0.463"
$ws.Cells.Item(12, 5).Value = "  +1.68%  "
$ws.Cells.Item(13, 5).Value = "  +15.21%  "
$ws.Cells.Item(14, 4).Value = "37.14"
$ws.Cells.Item(14, 5).Value = "  +1.09%  "
$ws.Cells.Item(15, 4).Value = "4.494.99"
$ws.Cells.Item(15, 5).Value = "  -0.02%  "
$ws.Cells.Item(16, 4).Value = "3.835.08"
$ws.Cells.Item(16, 5).Value = "  -0.73%  "
$ws.Cells.Item(17, 4).Value = "68.401.06"
$ws.Cells.Item(17, 5).Value = "  +0.74%  "
$ws.Cells.Item(18, 4).Value = "18.42"
$ws.Cells.Item(18, 5).Value = "  +1.89%  "
$ws.Cells.Item(19, 4).Value = "7.49"
$ws.Cells.Item(19, 5).Value = "  +2.02%  "
$ws.Cells.Item(20, 5).Value = "  +0.71%  "
$ws.Cells.Item(21, 5).Value = "  +0.15%  "
$ws.Cells.Item(22, 4).Value = "471.31"
$ws.Cells.Item(22, 5).Value = "  +1.82%  "
$ws.Cells.Item(23, 4).Value = "0.735"
$ws.Cells.Item(23, 5).Value = "  +1.01%  "
$ws.Cells.Item(24, 5).Value = "  -2.38%  "
$ws.Cells.Item(25, 4).Value = "83.70"
$ws.Cells.Item(25, 5).Value = "  +0.66%  "
$ws.Cells.Item(26, 4).Value = "2.29"
$ws.Cells.Item(26, 5).Value = "  +2.41%  "
$ws.Cells.Item(27, 4).Value = "12.21"
$ws.Cells.Item(27, 5).Value = "  +1.10%  "
$ws.Cells.Item(28, 4).Value = "10.56"
$ws.Cells.Item(28, 5).Value = "  +5.91%  "
$ws.Cells.Item(29, 5).Value = "  +0.05%  "
$ws.Cells.Item(30, 5).Value = "  -0.02%  "
$ws.Cells.Item(31, 4).Value = "3.999.63"
$ws.Cells.Item(31, 5).Value = "  -0.02%  "
$ws.Cells.Item(32, 4).Value = "7.78"
$ws.Cells.Item(32, 5).Value = "  +0.72%  "
$ws.Cells.Item(33, 5).Value = "  +0.05%  "
$ws.Cells.Item(34, 4).Value = "31.20"
$ws.Cells.Item(34, 5).Value = "  +0.75%  "
$ws.Cells.Item(35, 4).Value = "9.43"
$ws.Cells.Item(35, 5).Value = "  +1.36%  "
$ws.Cells.Item(36, 4).Value = "3.810.07"
$ws.Cells.Item(37, 4).Value = "3.93"
$ws.Cells.Item(37, 5).Value = "  +20.19%  "
$ws.Cells.Item(38, 4).Value = "0.106"
$ws.Cells.Item(38, 5).Value = "  +1.50%  "
$ws.Cells.Item(39, 4).Value = "5.99"
$ws.Cells.Item(39, 5).Value = "  +1.88%  "
$ws.Cells.Item(40, 5).Value = "  +0.28%  "
$ws.Cells.Item(41, 5).Value = "  +0.39%  "
$ws.Cells.Item(42, 4).Value = "0.999"
$ws.Cells.Item(42, 5).Value = "  -0.15%  "
$ws.Cells.Item(43, 4).Value = "0.321"
$ws.Cells.Item(43, 5).Value = "  +3.05%  "
$ws.Cells.Item(44, 5).Value = "  +0.87%  "
$ws.Cells.Item(46, 4).Value = "8.78"
$ws.Cells.Item(46, 5).Value = "  +3.53%  "
$ws.Cells.Item(47, 4).Value = "419.91"
$ws.Cells.Item(47, 5).Value = "  -1.61%  "
$ws.Cells.Item(48, 5).Value = "  +7.38%  "
$ws.Cells.Item(49, 4).Value = "46.77"
$ws.Cells.Item(49, 5).Value = "  -0.77%  "
$ws.Cells.Item(50, 4).Value = "142.29"
$ws.Cells.Item(50, 5).Value = "  -0.87%  "
$ws.Cells.Item(51, 5).Value = "  +2.03%  "
